# "Added TaskList page and Test Case"
# The diff shows that Sheet1 gains a header/data row: A1 = "admin", B1 = "manager"
# (stored as shared strings 0 and 1 respectively), which also grows the sheet's
# dimension from "A1" to "A1:B1" and populates xl/sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "admin"
$ws.Range("B1").Value = "manager"
